$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style helpers -------------------------------------------------------
# xlPasteFormats = -4122, xlRight = -4152
$xlPasteFormats = -4122
$xlRight = -4152

# 1) Highlight the category header cells in column A (rows 3,8,13,18,23,28)
#    Their new look = fill from B1 (light theme fill) + existing right alignment.
$headerCells = @("A3", "A8", "A13", "A18", "A23", "A28")
foreach ($addr in $headerCells) {
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).HorizontalAlignment = $xlRight
}
$excel.CutCopyMode = $false

# 2) Replace the rich-text "** Query trip #1" / " ** #2" / "** #3" labels
#    with plain labels, moving the "**" footnote marker into its own cell (E column).
$ws.Range("A28").Value2 = "Query trip #1"
$ws.Range("A29").Value2 = "#2"
$ws.Range("A30").Value2 = "#3"

# New "**" marker cells in column E, rows 28-31, styled like the existing
# "***" marker cells in column I (fontId4/red bold, bordered, right aligned).
$markerRows = @(28, 29, 30, 31)
foreach ($r in $markerRows) {
    $ws.Range("I28").Copy() | Out-Null
    $ws.Range("E$r").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("E$r").Value2 = "**"
}
$excel.CutCopyMode = $false

# New "***" marker cell in column I, row 31 (median row), matching I28:I30.
$ws.Range("I28").Copy() | Out-Null
$ws.Range("I31").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I31").Value2 = "***"
$excel.CutCopyMode = $false

# 3) Highlight the median totals on row 31 (J31 orange, K31 red) while keeping
#    their MEDIAN() formulas and number format intact.
$ws.Range("J31").Interior.Color = 49407   # FFFFC000 -> BGR 0x00C0FF
$ws.Range("K31").Interior.Color = 255     # FFFF0000 -> BGR 0x0000FF

# 4) Add a new row 32 with an italic note in J32 (same look as B37/B38 notes).
$ws.Range("B37").Copy() | Out-Null
$ws.Range("J32").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J32").Value2 = "Not fair to compare 128 vs 16668 records"
$excel.CutCopyMode = $false

# 5) Update the footnote text: the old "only returns max 128 rows" note is
#    replaced by a note about 1668 rows.
$ws.Range("B35").Value2 = "RavenDb returns 128 rows while SisoDb returns the full match of 1668 rows"

# 6) Restore the scroll position / active selection recorded in the workbook.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("J33").Select() | Out-Null

# 7) Match the saved window size recorded in the workbook view.
$excel.ActiveWindow.Width = 20730
$excel.ActiveWindow.Height = 11760
